$wb = $excel.ActiveWorkbook

# --- Sheet1 : update clearing cheque number, refresh selection ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A2").Value = "CLK0029898"
$ws1.Range("B30").Select() | Out-Null

# --- Sheet2 : update clearing cheque number, refresh selection ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A2").Value = "CLK0029896"
$ws2.Range("A2").Select() | Out-Null

# --- Sheet3 : update clearing cheque number, add a new data row, refresh selection ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A2").Value = "CLK0029894"

# New row 3: match the text number format used by B2 so "002" is kept as text
$ws3.Range("B3").NumberFormat = $ws3.Range("B2").NumberFormat
$ws3.Range("A3").Value = "CLK0029893"
$ws3.Range("B3").Value = "002"
$ws3.Range("C3").Value = 1

$ws3.Range("H19").Select() | Out-Null
